$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10138.667
$ws.Range("I51").Value = 9997.5
$ws.Range("J51").Value = 10151.5
$ws.Range("K51").Value = 9997.5
$ws.Range("L51").Value = 10151.5
$ws.Range("M51").Value = -9513.5
$ws.Range("N51").Value = -11119.5

$ws.Range("H76").Value = 3353.3
$ws.Range("I76").Value = 2878
$ws.Range("J76").Value = 3828.6
$ws.Range("K76").Value = 2878
$ws.Range("L76").Value = 3828.6
$ws.Range("M76").Value = -2563
$ws.Range("N76").Value = -4458.6

$ws.Range("H79").Value = 3353.3
$ws.Range("I79").Value = 2878
$ws.Range("J79").Value = 3828.6
$ws.Range("K79").Value = 2878
$ws.Range("L79").Value = 3828.6
$ws.Range("M79").Value = -1786
$ws.Range("N79").Value = -6012.6

$ws.Range("H112").Value = 4107.08
$ws.Range("J112").Value = 4107.375
$ws.Range("L112").Value = 12322.125
$ws.Range("N112").Value = -14538.125

$ws.Range("H116").Value = 94452750
$ws.Range("I116").Value = 113341736
$ws.Range("K116").Value = 113341736
$ws.Range("M116").Value = -113338294

$ws.Range("H136").Value = 150000
$ws.Range("J136").Value = 150000
$ws.Range("L136").Value = 150000
$ws.Range("N136").Value = -160200

$ws.Range("H138").Value = 5955.2856
$ws.Range("I138").Value = 1758.8235
$ws.Range("J138").Value = 7144.283
$ws.Range("K138").Value = 5276.470499999999
$ws.Range("L138").Value = 21432.849
$ws.Range("M138").Value = -136.4704999999994
$ws.Range("N138").Value = -31712.849

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H32").Value = 18096.877
$ws.Range("I32").Value = 17130.658
$ws.Range("K32").Value = 17130.658
$ws.Range("M32").Value = -16843.658

$ws.Range("H37").Value = 31124.75
$ws.Range("I37").Value = 31124.75
$ws.Range("K37").Value = 31124.75
$ws.Range("M37").Value = -30851.75

$ws.Range("H74").Value = 35716460
$ws.Range("I74").Value = 35716460
$ws.Range("K74").Value = 35716460
$ws.Range("M74").Value = -35715586

$ws.Range("H77").Value = 35716460
$ws.Range("I77").Value = 35716460
$ws.Range("K77").Value = 178582300
$ws.Range("M77").Value = -178577932

$ws.Range("H132").Value = 1825.0889
$ws.Range("I132").Value = 1866
$ws.Range("J132").Value = 945.5
$ws.Range("K132").Value = 5598
$ws.Range("L132").Value = 2836.5
$ws.Range("M132").Value = -3068
$ws.Range("N132").Value = -7896.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 99998.5
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()

$ws.Range("H107").Value = 2237.805
$ws.Range("I107").Value = 2143.9644
$ws.Range("K107").Value = 2143.9644
$ws.Range("M107").Value = -223.9643999999998

$ws.Range("H134").Value = 2122.3784
$ws.Range("I134").Value = 1593.6296
$ws.Range("J134").Value = 3550
$ws.Range("K134").Value = 4780.8888
$ws.Range("L134").Value = 10650
$ws.Range("M134").Value = -2245.8888
$ws.Range("N134").Value = -15720

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 11998.6
$ws.Range("I6").Value = 11998.6
$ws.Range("K6").Value = 11998.6
$ws.Range("M6").Value = -11885.6

$ws.Range("H7").Value = 232.06667
$ws.Range("I7").Value = 155.6
$ws.Range("J7").Value = 385
$ws.Range("K7").Value = 155.6
$ws.Range("L7").Value = 385
$ws.Range("M7").Value = -42.59999999999999
$ws.Range("N7").Value = -611

$ws.Range("H58").Value = 2533.8948
$ws.Range("I58").Value = 2352.625
$ws.Range("K58").Value = 2352.625
$ws.Range("M58").Value = -2149.625

$ws.Range("H62").Value = 49537.5
$ws.Range("I62").Value = 24718.8
$ws.Range("K62").Value = 24718.8
$ws.Range("M62").Value = -24094.8

$ws.Range("H65").Value = 49537.5
$ws.Range("I65").Value = 24718.8
$ws.Range("K65").Value = 123594
$ws.Range("M65").Value = -120474

$ws.Range("H107").Value = 815.7
$ws.Range("J107").Value = 886.4286
$ws.Range("L107").Value = 886.4286
$ws.Range("N107").Value = -4726.4286

$ws.Range("H108").Value = 5391.6665
$ws.Range("J108").Value = 5391.6665
$ws.Range("L108").Value = 5391.6665
$ws.Range("N108").Value = -13071.6665

$ws.Range("H136").Value = 2533.8948
$ws.Range("I136").Value = 2352.625
$ws.Range("K136").Value = 7057.875
$ws.Range("M136").Value = -4507.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 276.53845
$ws.Range("J38").Value = 585
$ws.Range("L38").Value = 1755
$ws.Range("N38").Value = -2449

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H81").Value = 4435.1577
$ws.Range("J81").Value = 4746.25
$ws.Range("L81").Value = 14238.75
$ws.Range("N81").Value = -16484.75

$ws.Range("H84").Value = 4435.1577
$ws.Range("J84").Value = 4746.25
$ws.Range("L84").Value = 42716.25
$ws.Range("N84").Value = -53948.25

$ws.Range("H131").Value = 2960.6155
$ws.Range("I131").Value = 2734.7058
$ws.Range("J131").Value = 3070.3428
$ws.Range("K131").Value = 8204.117400000001
$ws.Range("L131").Value = 9211.028399999999
$ws.Range("M131").Value = -3164.117400000001
$ws.Range("N131").Value = -19291.0284

$ws.Range("H134").Value = 12356.066
$ws.Range("I134").Value = 9641.615
$ws.Range("K134").Value = 28924.845
$ws.Range("M134").Value = -23854.845

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 316731.72
$ws.Range("J122").Value = 6246
$ws.Range("L122").Value = 18738
$ws.Range("N122").Value = -23638

$ws.Range("H132").Value = 2851.0667
$ws.Range("I132").Value = 1580.5
$ws.Range("J132").Value = 7933.3335
$ws.Range("K132").Value = 4741.5
$ws.Range("L132").Value = 23800.0005
$ws.Range("M132").Value = -2211.5
$ws.Range("N132").Value = -28860.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 200000
$ws.Range("I11").Value = 200000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 200000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -199860
$ws.Range("N11").ClearContents()

$ws.Range("H16").Value = 78759.234
$ws.Range("I16").Value = 85072.5
$ws.Range("K16").Value = 85072.5
$ws.Range("M16").Value = -84902.5

$ws.Range("H20").Value = 15699.286
$ws.Range("I20").Value = 15699.286
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 15699.286
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -15473.286
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 862.06665
$ws.Range("I22").Value = 792.55554
$ws.Range("J22").Value = 966.3333
$ws.Range("K22").Value = 792.55554
$ws.Range("L22").Value = 966.3333
$ws.Range("M22").Value = -497.55554
$ws.Range("N22").Value = -1556.3333

$ws.Range("H26").Value = 25900
$ws.Range("I26").Value = 4800
$ws.Range("J26").Value = 47000
$ws.Range("K26").Value = 4800
$ws.Range("L26").Value = 47000
$ws.Range("M26").Value = -4505
$ws.Range("N26").Value = -47590

$ws.Range("H27").Value = 862.06665
$ws.Range("I27").Value = 792.55554
$ws.Range("J27").Value = 966.3333
$ws.Range("K27").Value = 792.55554
$ws.Range("L27").Value = 966.3333
$ws.Range("M27").Value = -685.55554
$ws.Range("N27").Value = -1180.3333

$ws.Range("H40").Value = 19930722
$ws.Range("I40").Value = 7814665.5
$ws.Range("K40").Value = 7814665.5
$ws.Range("M40").Value = -7814529.5

$ws.Range("H68").Value = 5050.385
$ws.Range("I68").Value = 3471.25
$ws.Range("J68").Value = 7577
$ws.Range("K68").Value = 3471.25
$ws.Range("L68").Value = 7577
$ws.Range("M68").Value = -2722.25
$ws.Range("N68").Value = -9075

$ws.Range("H71").Value = 5050.385
$ws.Range("I71").Value = 3471.25
$ws.Range("J71").Value = 7577
$ws.Range("K71").Value = 17356.25
$ws.Range("L71").Value = 37885
$ws.Range("M71").Value = -13612.25
$ws.Range("N71").Value = -45373

$ws.Range("H136").Value = 4088.8262
$ws.Range("I136").Value = 2802.4473
$ws.Range("J136").Value = 5665.6772
$ws.Range("K136").Value = 8407.341899999999
$ws.Range("L136").Value = 16997.0316
$ws.Range("M136").Value = -5857.341899999999
$ws.Range("N136").Value = -22097.0316

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 100300
$ws.Range("J80").Value = 100300
$ws.Range("L80").Value = 100300
$ws.Range("N80").Value = -102296

$ws.Range("H81").Value = 12511.714
$ws.Range("I81").Value = 11264.5
$ws.Range("K81").Value = 22529
$ws.Range("M81").Value = -21468

$ws.Range("H83").Value = 100300
$ws.Range("J83").Value = 100300
$ws.Range("L83").Value = 300900
$ws.Range("N83").Value = -310884

$ws.Range("H84").Value = 12511.714
$ws.Range("I84").Value = 11264.5
$ws.Range("K84").Value = 112645
$ws.Range("M84").Value = -107341

$ws.Range("H126").Value = 2045.2142
$ws.Range("I126").Value = 1827.4286
$ws.Range("K126").Value = 5482.2858
$ws.Range("M126").Value = -3012.2858

$ws.Range("H132").Value = 5653925
$ws.Range("I132").Value = 13334331
$ws.Range("K132").Value = 40002993
$ws.Range("M132").Value = -40000463

$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws.Range("H136").Value = 3587.6724
$ws.Range("I136").Value = 2375.6667
$ws.Range("J136").Value = 5570.9546
$ws.Range("K136").Value = 7127.000100000001
$ws.Range("L136").Value = 16712.8638
$ws.Range("M136").Value = -4577.000100000001
$ws.Range("N136").Value = -21812.8638
